$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "35÷3=11, 2" "33÷5=6, 3"
Replace-Text "56÷9=6, 2" "92÷5=18, 2"
Replace-Text "50÷5=10, 0" "72÷4=18, 0"
Replace-Text "20÷8=2, 4" "24÷8=3, 0"
Replace-Text "68÷8=8, 4" "66÷2=33, 0"
Replace-Text "75÷5=15, 0" "57÷4=14, 1"
Replace-Text "26÷6=4, 2" "75÷4=18, 3"
Replace-Text "99÷8=12, 3" "89÷4=22, 1"
Replace-Text "94÷5=18, 4" "13÷8=1, 5"
Replace-Text "95÷8=11, 7" "68÷4=17, 0"
Replace-Text "48÷6=8, 0" "25÷6=4, 1"
Replace-Text "93÷7=13, 2" "71÷8=8, 7"
Replace-Text "36÷2=18, 0" "16÷3=5, 1"
Replace-Text "99÷9=11, 0" "73÷2=36, 1"
Replace-Text "19÷4=4, 3" "58÷5=11, 3"
Replace-Text "96÷8=12, 0" "64÷7=9, 1"
Replace-Text "40÷8=5, 0" "66÷8=8, 2"
Replace-Text "95÷2=47, 1" "51÷7=7, 2"
Replace-Text "24÷4=6, 0" "75÷3=25, 0"
Replace-Text "52÷9=5, 7" "14÷6=2, 2"
Replace-Text "84÷3=28, 0" "42÷4=10, 2"
Replace-Text "97÷9=10, 7" "54÷7=7, 5"
Replace-Text "83÷7=11, 6" "53÷6=8, 5"
Replace-Text "18÷5=3, 3" "71÷7=10, 1"
Replace-Text "19÷9=2, 1" "64÷3=21, 1"
